$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G4").Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G5").Value = 'nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G6").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G7").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G8").Value = 'manar.montaser@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, backup@backdoor.com, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G12").Value = 'salma.elgendy.std@med.asu.edu.eg, System'
$ws.Range("G13").Value = 'mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G20").Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G24").Value = 'maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G25").Value = 'Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G27").Value = 'Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G30").Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G31").Value = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G32").Value = 'Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G33").Value = 'nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, hananragab@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G34").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G35").Value = 'Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G36").Value = 'manar.montaser@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, backup@backdoor.com, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G40").Value = 'salma.elgendy.std@med.asu.edu.eg, System'
$ws.Range("G41").Value = 'mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G48").Value = 'nourhan.mostafa@med.asu.edu.eg, user@user.com'
$ws.Range("G52").Value = 'maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G53").Value = 'Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G55").Value = 'Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
